# Insert one new data row at sheet row 186 (pushing the existing rows
# 186-333 down to 187-334, i.e. a new "weekly" Ciboulette price record is
# slotted in chronologically and everything after it shifts down by one
# row). This grows the used range from A1:R333 to A1:R334.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pushes old row 186 (and everything below it) down to row 187, leaving a
# blank row 186 ready to be filled in with the new record.
$ws.Rows.Item(186).Insert()

# Populate the newly-inserted row 186 with the new weekly record. All
# columns other than "Fecha" (D) and "Volumen" (J) repeat the same
# constant template used by every other row in this data set.
$ws.Cells.Item(186, 1).Value = 3
$ws.Cells.Item(186, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(186, 3).Value = "Coquimbo"
$ws.Cells.Item(186, 4).Value = 44741
$ws.Cells.Item(186, 5).Value = 5
$ws.Cells.Item(186, 6).Value = 100112039
$ws.Cells.Item(186, 7).Value = "Ciboulette"
$ws.Cells.Item(186, 8).Value = "Sin especificar"
$ws.Cells.Item(186, 9).Value = "Primera"
$ws.Cells.Item(186, 10).Value = 120
$ws.Cells.Item(186, 11).Value = 1500
$ws.Cells.Item(186, 12).Value = 1500
$ws.Cells.Item(186, 13).Value = 1500
$ws.Cells.Item(186, 14).Value = "`$/docena de atados"
$ws.Cells.Item(186, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(186, 16).Value = 500
$ws.Cells.Item(186, 17).Value = 3
$ws.Cells.Item(186, 18).Value = "Hortaliza"
